$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: "Correct %"  = Correct / (Correct + Multiple), per folder ---
# --- Row 13: "Multiple %" = Multiple / (Multiple + Correct), per folder ---
# Copy formatting (fill/border/style) from row 11 down onto the two new rows
$ws.Range("A11:AG11").Copy() | Out-Null
$ws.Range("A12:A13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A12").Value = "Correct %"
$ws.Range("B12").Formula = "=B4/(B4+B6)"
$ws.Range("C12:AG12").Formula = "=C4/(C4+C6)"

$ws.Range("A13").Value = "Multiple %"
$ws.Range("B13").Formula = "=B6/(B6+B4)"
$ws.Range("C13:AG13").Formula = "=C6/(C6+C4)"

# --- Rows 19-22: per-file / per-folder correct vs indefinite summary rows ---
# Copy formatting from row 18 (single-column summary row) onto the new rows
$ws.Range("A18:B18").Copy() | Out-Null
$ws.Range("A19:A22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A19").Value = "Correct/(correct+indef) per file"
$ws.Range("B19").Formula = "=SUM(B4:AG4)/(SUM(B4:AG4)+SUM(B6:AG6))"

$ws.Range("A20").Value = "Indef/(correct+indef) per file"
$ws.Range("B20").Formula = "=1-B19"

$ws.Range("A21").Value = "Correct/(correct+indef) per folder"
$ws.Range("B21").Formula = "=AVERAGE(B12:AG12)"

$ws.Range("A22").Value = "Indef/(correct+indef) per folder"
$ws.Range("B22").Formula = "=AVERAGE(B13:AG13)"

# --- View changes: scrolled/selected state moved down a bit ---
# (topLeftCell goes from B1 -> A3, selection goes from E20 -> E21)
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("E21").Select() | Out-Null
